# Fruta / hortaliza, semanal
# Insert the new weekly price-report row for
# "Feria Lagunitas de Puerto Montt - Mango" at row 409, pushing the
# existing rows 409:425 down to 410:426.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 409 (shifts 409:425 -> 410:426, carrying
# over cell formatting/styles such as the date format on column D).
$ws.Rows.Item(409).Insert()

# Populate the new row 409 with this week's data.
$ws.Cells.Item(409, 1).Value = 4
$ws.Cells.Item(409, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(409, 3).Value = "Los Lagos"
$ws.Cells.Item(409, 4).Value = 45147
$ws.Cells.Item(409, 5).Value = 10
$ws.Cells.Item(409, 6).Value = "Fruta"
$ws.Cells.Item(409, 7).Value = 100108
$ws.Cells.Item(409, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(409, 9).Value = 100108002
$ws.Cells.Item(409, 10).Value = "Mango"
$ws.Cells.Item(409, 11).Value = "Sin especificar"
$ws.Cells.Item(409, 12).Value = "Primera"
$ws.Cells.Item(409, 13).Value = 40
$ws.Cells.Item(409, 14).Value = 10000
$ws.Cells.Item(409, 15).Value = 10000
$ws.Cells.Item(409, 16).Value = 10000
$ws.Cells.Item(409, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(409, 18).Value = "Perú"
$ws.Cells.Item(409, 19).Value = 2500
$ws.Cells.Item(409, 20).Value = 4
